# Applies updated leve-profit market data for each job sheet,
# refreshing price/profit columns (H:N) on the affected rows.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 150
$ws.Range("I33").Value = 150
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 150
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 79
$ws.Range("N33").ClearContents()

$ws.Range("H138").Value = 2814.9534
$ws.Range("I138").Value = 2426.5518
$ws.Range("J138").Value = 3619.5
$ws.Range("K138").Value = 7279.655400000001
$ws.Range("L138").Value = 10858.5
$ws.Range("M138").Value = -2139.655400000001
$ws.Range("N138").Value = -21138.5

$ws.Range("H141").Value = 905873.5
$ws.Range("I141").Value = 1168435
$ws.Range("K141").Value = 3505305
$ws.Range("M141").Value = -3500125


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3028.2124
$ws.Range("I32").Value = 2369.6567
$ws.Range("J32").Value = 6422.3076
$ws.Range("K32").Value = 2369.6567
$ws.Range("L32").Value = 6422.3076
$ws.Range("M32").Value = -2082.6567
$ws.Range("N32").Value = -6996.3076

$ws.Range("H74").Value = 1481.4615
$ws.Range("I74").Value = 1230.4706
$ws.Range("J74").Value = 1955.5555
$ws.Range("K74").Value = 1230.4706
$ws.Range("L74").Value = 1955.5555
$ws.Range("M74").Value = -356.4706000000001
$ws.Range("N74").Value = -3703.5555

$ws.Range("H77").Value = 1481.4615
$ws.Range("I77").Value = 1230.4706
$ws.Range("J77").Value = 1955.5555
$ws.Range("K77").Value = 6152.353000000001
$ws.Range("L77").Value = 9777.7775
$ws.Range("M77").Value = -1784.353000000001
$ws.Range("N77").Value = -18513.7775

$ws.Range("H132").Value = 2654.611
$ws.Range("I132").Value = 2665.8333
$ws.Range("J132").Value = 2649
$ws.Range("K132").Value = 7997.499899999999
$ws.Range("L132").Value = 7947
$ws.Range("M132").Value = -5467.499899999999
$ws.Range("N132").Value = -13007


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 419.64865
$ws.Range("I94").Value = 337.2353
$ws.Range("J94").Value = 1353.6666
$ws.Range("K94").Value = 337.2353
$ws.Range("L94").Value = 1353.6666
$ws.Range("M94").Value = 113.7647
$ws.Range("N94").Value = -2255.6666


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 850.6
$ws.Range("I16").Value = 813.625
$ws.Range("J16").Value = 998.5
$ws.Range("K16").Value = 813.625
$ws.Range("L16").Value = 998.5
$ws.Range("M16").Value = -526.625
$ws.Range("N16").Value = -1572.5

$ws.Range("H31").Value = 2028.0555
$ws.Range("I31").Value = 1362.375
$ws.Range("J31").Value = 2560.6
$ws.Range("K31").Value = 1362.375
$ws.Range("L31").Value = 2560.6
$ws.Range("M31").Value = -1067.375
$ws.Range("N31").Value = -3150.6

$ws.Range("H34").Value = 2028.0555
$ws.Range("I34").Value = 1362.375
$ws.Range("J34").Value = 2560.6
$ws.Range("K34").Value = 1362.375
$ws.Range("L34").Value = 2560.6
$ws.Range("M34").Value = -1160.375
$ws.Range("N34").Value = -2964.6

$ws.Range("H44").Value = 10280
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 10280
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 10280
$ws.Range("N44").Value = -11164
$ws.Range("M44").ClearContents()

$ws.Range("H50").Value = 19358.5
$ws.Range("J50").Value = 19358.5
$ws.Range("L50").Value = 19358.5
$ws.Range("N50").Value = -20608.5

$ws.Range("H103").Value = 18500
$ws.Range("I103").Value = 17000
$ws.Range("J103").Value = 20000
$ws.Range("K103").Value = 17000
$ws.Range("L103").Value = 20000
$ws.Range("M103").Value = -15828
$ws.Range("N103").Value = -22344

$ws.Range("H113").Value = 850.6
$ws.Range("I113").Value = 813.625
$ws.Range("J113").Value = 998.5
$ws.Range("K113").Value = 813.625
$ws.Range("L113").Value = 998.5
$ws.Range("M113").Value = 1356.375
$ws.Range("N113").Value = -5338.5

$ws.Range("H122").Value = 2975.15
$ws.Range("I122").Value = 1862.875
$ws.Range("K122").Value = 5588.625
$ws.Range("M122").Value = -3138.625

$ws.Range("H132").Value = 3746.5454
$ws.Range("I132").Value = 1450
$ws.Range("J132").Value = 5058.857
$ws.Range("K132").Value = 4350
$ws.Range("L132").Value = 15176.571
$ws.Range("M132").Value = -1820
$ws.Range("N132").Value = -20236.571

$ws.Range("H141").Value = 42928.91
$ws.Range("I141").Value = 10000
$ws.Range("J141").Value = 70369.664
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 70369.664
$ws.Range("M141").Value = -4820
$ws.Range("N141").Value = -80729.664


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 1200
$ws.Range("N12").Value = -1546
$ws.Range("M12").ClearContents()

$ws.Range("H122").Value = 1012.1
$ws.Range("J122").Value = 1296.8334
$ws.Range("L122").Value = 11671.5006
$ws.Range("N122").Value = -16571.5006

$ws.Range("H129").Value = 28822.5
$ws.Range("I129").Value = 534.375
$ws.Range("J129").Value = 41395
$ws.Range("K129").Value = 1603.125
$ws.Range("L129").Value = 124185
$ws.Range("M129").Value = 3396.875
$ws.Range("N129").Value = -134185

$ws.Range("H131").Value = 790.67
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 799.6598
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2398.9794
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12478.9794

$ws.Range("H137").Value = 3630.7778
$ws.Range("J137").Value = 3969.6
$ws.Range("L137").Value = 11908.8
$ws.Range("N137").Value = -22108.8


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 775.5
$ws.Range("I107").Value = 250
$ws.Range("J107").Value = 1301
$ws.Range("K107").Value = 250
$ws.Range("L107").Value = 1301
$ws.Range("M107").Value = 1670
$ws.Range("N107").Value = -5141

$ws.Range("H132").Value = 1925391.4
$ws.Range("I132").Value = 4275328
$ws.Range("J132").Value = 2716.182
$ws.Range("K132").Value = 12825984
$ws.Range("L132").Value = 8148.545999999999
$ws.Range("M132").Value = -12823454
$ws.Range("N132").Value = -13208.546

$ws.Range("H136").Value = 11189.429
$ws.Range("J136").Value = 11189.429
$ws.Range("L136").Value = 33568.287
$ws.Range("N136").Value = -38668.287


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1999.3077
$ws.Range("I61").Value = 1899.3
$ws.Range("K61").Value = 1899.3
$ws.Range("M61").Value = -1697.3

$ws.Range("H113").Value = 1999.3077
$ws.Range("I113").Value = 1899.3
$ws.Range("K113").Value = 1899.3
$ws.Range("M113").Value = 270.7

$ws.Range("H132").Value = 1762.9445
$ws.Range("I132").Value = 1580.2142
$ws.Range("J132").Value = 1879.2273
$ws.Range("K132").Value = 4740.642599999999
$ws.Range("L132").Value = 5637.6819
$ws.Range("M132").Value = -2210.642599999999
$ws.Range("N132").Value = -10697.6819


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 13493.333
$ws.Range("J29").Value = 18995
$ws.Range("L29").Value = 18995
$ws.Range("N29").Value = -19575

$ws.Range("H81").Value = 1730.7778
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 1730.7778
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 30000
$ws.Range("N84").Value = -40608

$ws.Range("H107").Value = 627.44446
$ws.Range("J107").Value = 921.2
$ws.Range("L107").Value = 2763.6
$ws.Range("N107").Value = -6603.6

$ws.Range("H126").Value = 13409.272
$ws.Range("I126").Value = 16571.715
$ws.Range("K126").Value = 49715.145
$ws.Range("M126").Value = -47245.145

$ws.Range("H132").Value = 1389.2222
$ws.Range("I132").Value = 1125.4375
$ws.Range("J132").Value = 3499.5
$ws.Range("K132").Value = 3376.3125
$ws.Range("L132").Value = 10498.5
$ws.Range("M132").Value = -846.3125
$ws.Range("N132").Value = -15558.5

$ws.Range("H136").Value = 29243064
$ws.Range("I136").Value = 39685988
$ws.Range("J136").Value = 2880
$ws.Range("K136").Value = 119057964
$ws.Range("L136").Value = 8640
$ws.Range("M136").Value = -119055414
$ws.Range("N136").Value = -13740
